# Apply "average with safety stocks" update.
$wb = $excel.ActiveWorkbook

# --- Productdata sheet: rescale Leadtimes(D)/InventoryCosts(F)/SetupCosts(I) columns ---
$productData = $wb.Worksheets.Item("Productdata")

$productDataUpdates = @{
    "D2" = 0.00448
    "F2" = 0.008959999999999999
    "I2" = 0.08959999999999999

    "D3" = 0.00496
    "F3" = 0.00992
    "I3" = 0.0992

    "D4" = 0.004920000000000001
    "F4" = 0.009840000000000002
    "I4" = 0.09840000000000002

    "D5" = 0.00444
    "F5" = 0.008880000000000001
    "I5" = 0.0888

    "D6" = 0.00048
    "F6" = 0.00096
    "I6" = 0.009600000000000001

    "D7" = 0.00048
    "F7" = 0.00096
    "I7" = 0.009600000000000001

    "D8" = 0.0004400000000000001
    "F8" = 0.0008800000000000001
    "I8" = 0.008800000000000002

    "D9" = 0.00004
    "F9" = 0.00008000000000000001
    "I9" = 0.0008

    "D10" = 0.00004
    "F10" = 0.00008000000000000001
    "I10" = 0.0008

    "D11" = 0.00004
    "F11" = 0.00008000000000000001
    "I11" = 0.0008
}

foreach ($addr in $productDataUpdates.Keys) {
    $productData.Range($addr).Value = $productDataUpdates[$addr]
}

# --- ForcastedStandardDeviation sheet: zero out rows 9-11, columns B-E ---
$forecastStdDev = $wb.Worksheets.Item("ForcastedStandardDeviation")

foreach ($row in 9..11) {
    foreach ($col in @("B", "C", "D", "E")) {
        $forecastStdDev.Range("$col$row").Value = 0
    }
}
